# Automatische test-sync: 2025-08-05 19:38:50
# Appends Testmail #8 ("Check jij even of dit nog geleverd kan worden?")
# as a new row to the "Logs" sheet, bumps the "Overig" tally on the
# "Dashboard" sheet, and extends the conditional-formatting ranges on
# "Logs" to cover the newly added row.

function ToBgr($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $b * 65536 + $g * 256 + $r
}

$wb = $excel.ActiveWorkbook

# --- Logs sheet: append new row 49 with the new test mail entry ---
$logs = $wb.Worksheets.Item("Logs")
$newRow = 49

$logs.Cells.Item($newRow, 1).Value = "Check jij even of dit nog geleverd kan worden?"
$logs.Cells.Item($newRow, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item($newRow, 3).Value = "Testmail #8: Check jij even of dit nog geleverd kan worden?"
$logs.Cells.Item($newRow, 4).Value = "Overig"
$logs.Cells.Item($newRow, 5).Value = "Beste klant,`nDank u voor uw e-mail. Om uw vraag over de levering te beantwoorden, hebben we wat meer informatie nodig. Kunt u ons alstublieft het volgende verstrekken:`n- Het specifieke product dat u wilt bestellen`n- Het afleveradres`n- Gewenste leverdatum`nMet deze informatie kunnen we controleren of het product nog steeds geleverd kan worden. We kijken uit naar uw antwoord.`nMet vriendelijke groet,`n[Bedrijfsnaam] E-mailassistent"
$logs.Cells.Item($newRow, 6).Value = "2025-08-05 19:37:52"
$logs.Cells.Item($newRow, 7).Value = "Ja"
$logs.Cells.Item($newRow, 8).Value = "Nee"
$logs.Cells.Item($newRow, 9).Value = "Ja"
$logs.Cells.Item($newRow, 10).Value = "Nee"

# Undo any automatic row-height growth triggered by the multi-line answer
# text in column E so the row keeps the sheet's default (non-custom) height.
$logs.Rows.Item($newRow).AutoFit()

# --- Extend the existing conditional formats so they keep covering the
#     whole data range (...2:...48 -> ...2:...49), preserving the original
#     rule order/colors/priorities. ---
$logs.Range("D2:D48").FormatConditions.Delete()
$logs.Range("G2:G48").FormatConditions.Delete()
$logs.Range("H2:H48").FormatConditions.Delete()
$logs.Range("I2:I48").FormatConditions.Delete()
$logs.Range("J2:J48").FormatConditions.Delete()

$fcs = $logs.Range("D2:D49").FormatConditions
$fc = $fcs.Add(1, 3, '="Klacht"');             $fc.Interior.Color = ToBgr("FFC7CE")
$fc = $fcs.Add(1, 3, '="Bestelling"');         $fc.Interior.Color = ToBgr("FFEB9C")
$fc = $fcs.Add(1, 3, '="Informatieaanvraag"'); $fc.Interior.Color = ToBgr("C6EFCE")
$fc = $fcs.Add(1, 3, '="Afmelding"');          $fc.Interior.Color = ToBgr("BDD7EE")
$fc = $fcs.Add(1, 3, '="Retour"');             $fc.Interior.Color = ToBgr("FFC7CE")
$fc = $fcs.Add(1, 3, '="Overig"');             $fc.Interior.Color = ToBgr("D9D9D9")

$fcs = $logs.Range("G2:G49").FormatConditions
$fc = $fcs.Add(1, 3, '="Ja"');  $fc.Interior.Color = ToBgr("C6EFCE")
$fc = $fcs.Add(1, 3, '="Nee"'); $fc.Interior.Color = ToBgr("FFC7CE")

$fcs = $logs.Range("H2:H49").FormatConditions
$fc = $fcs.Add(1, 3, '="Ja"'); $fc.Interior.Color = ToBgr("FFF2CC")

$fcs = $logs.Range("I2:I49").FormatConditions
$fc = $fcs.Add(1, 3, '="Ja"'); $fc.Interior.Color = ToBgr("C9DAF8")

$fcs = $logs.Range("J2:J49").FormatConditions
$fc = $fcs.Add(1, 3, '="Ja"'); $fc.Interior.Color = ToBgr("D9D2E9")

# --- Dashboard sheet: bump the "Overig" category count from 6 to 7 ---
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Cells.Item(3, 2).Value = 7
